# Generate Report for Handoff
# Updates the "b.md" rows across the Overview, zh-cn, and de-de sheets to
# reflect that the file has been newly handed off (status + new handoff
# file name + new handoff datetime), replacing the old "Handed back" state.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: b.md row (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-27-19 10:27:06"

# --- zh-cn sheet: b.md row (row 3) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Item(3).Delete()
$wsZh.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-19 10:27:03"
$wsZh.Range("C3").Value = "Ready for handoff"
$zhTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/175a29edb7578ded6ee6dd76f8d64d2ec98eb740/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhTarget, "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf")

# --- de-de sheet: b.md row (row 3) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Item(3).Delete()
$wsDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-19 10:27:06"
$wsDe.Range("C3").Value = "Ready for handoff"
$deTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b19aabe2acd25c7e42c915c3286866c2e14fb0b3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deTarget, "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf")
